$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original (default) style, then force Text format so Excel
# does not reinterpret numeric-looking strings as numbers (which would
# drop significant trailing zeros / change representation).
$styleD = $ws.Range("D2").Style
$styleE = $ws.Range("E2").Style
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "24.609.33"
$ws.Range("E2").Value = "  +3.54%  "
$ws.Range("D3").Value = "1.696.56"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "318.54"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "0.3956"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("D8").Value = "0.4020"
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("D9").Value = "1.540"
$ws.Range("E9").Value = "  +8.33%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").Value = "53.48"
$ws.Range("E11").Value = "  +7.59%  "
$ws.Range("D12").Value = "0.08797"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").Value = "7.270"
$ws.Range("E13").Value = "  +8.70%  "
$ws.Range("D14").Value = "23.29"
$ws.Range("E14").Value = "  +2.35%  "
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "7.635"
$ws.Range("E16").Value = "  +5.45%  "
$ws.Range("D17").Value = "1.697.15"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "101.15"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").Value = "0.07008"
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").Value = "19.72"
$ws.Range("E20").Value = "  +3.12%  "
$ws.Range("D21").Value = "6.903"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("D24").Value = "24.624.10"
$ws.Range("E24").Value = "  +3.72%  "
$ws.Range("D25").Value = "3.071"
$ws.Range("E25").Value = "  +7.89%  "
$ws.Range("D26").Value = "2.339"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").Value = "22.39"
$ws.Range("E27").Value = "  +2.80%  "
$ws.Range("D28").Value = "159.81"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").Value = "5.246"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").Value = "134.57"
$ws.Range("E30").Value = "  +3.81%  "
$ws.Range("D31").Value = "7.467"
$ws.Range("E31").Value = "  +15.15%  "
$ws.Range("D32").Value = "1.107"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("D33").Value = "1.883.57"
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").Value = "7.524"
$ws.Range("E34").Value = "  +15.59%  "
$ws.Range("D35").Value = "0.08546"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").Value = "11.49"
$ws.Range("E36").Value = "  +10.87%  "
$ws.Range("D37").Value = "1.973"
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").Value = "0.2744"
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("D39").Value = "14.61"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").Value = "0.02776"
$ws.Range("E40").Value = "  +9.67%  "
$ws.Range("D41").Value = "0.09042"
$ws.Range("E41").Value = "  +2.92%  "
$ws.Range("D42").Value = "1.466"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").Value = "0.7719"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("D44").Value = "0.7233"
$ws.Range("E44").Value = "  +2.72%  "
$ws.Range("E45").Value = "  +3.82%  "
$ws.Range("D46").Value = "2.523"
$ws.Range("E46").Value = "  +5.54%  "
$ws.Range("D47").Value = "4.240"
$ws.Range("E47").Value = "  +3.64%  "
$ws.Range("D48").Value = "1.355"
$ws.Range("E48").Value = "  +13.68%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "141.24"
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("D51").Value = "0.08045"
$ws.Range("E51").Value = "  +3.72%  "

# Restore original style/number format so only the cell contents change.
$ws.Range("D2:D51").Style = $styleD
$ws.Range("E2:E51").Style = $styleE
